$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 12705
$wsExpo.Range("F10").Value = 12605
$wsExpo.Range("F14").Value = 5877

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 12705
$wsAll.Range("F11").Value = 12605
$wsAll.Range("F15").Value = 5877
